# repull data, push all data, mean calculation
# Update column F (dSF) values for several rows to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -4
$ws.Range("F6").Value = -3
$ws.Range("F12").Value = 2
$ws.Range("F15").Value = -4
$ws.Range("F17").Value = -3
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = -4
$ws.Range("F20").Value = 0
